# Update cryptos list values (Price / Volume(1h)) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.783.73"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "2.079.59"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.39"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.69"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  +2.78%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.04"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.775"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.35"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "2.046.58"
$ws.Range("E16").Value = "  -3.16%  "
$ws.Range("D17").Value = "37.715.64"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.60"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.40"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("E24").Value = "  -5.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.58"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.25"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.136"
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.42"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.50"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.82"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.39"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.37"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0978"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.52"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.64"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("D44").Value = "1.440.65"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.21"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "2.269.56"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.53"
$ws.Range("E51").Value = "  -0.15%  "
